# working on the new soma
# 1) A383:A387 were stored as text ("20081600" etc.) - convert them to real numbers.
# 2) Append new data rows 392:398 (date in col A as text, issues/maturities as numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix A383:A387: re-store the textual date codes as genuine numbers ---
$fixRows = 383..387
foreach ($r in $fixRows) {
    $cell = $ws.Cells.Item($r, 1)
    $numVal = $cell.Value2
    $cell.Value = $numVal
}

# --- Append rows 392:398 ---
$newData = @(
    @(392, "20082500", 120000000000, 110000000000),
    @(393, "20082600", 0, 0),
    @(394, "20082700", 165000000000, 162000000000),
    @(395, "20082800", 22000000000, 0),
    @(396, "20082900", 0, 0),
    @(397, "20083000", 0, 0),
    @(398, "20083100", 180000000000, 100000000000)
)

foreach ($row in $newData) {
    $r = $row[0]

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[1]
    $cellA.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

$wb.Save()
